# Add files for routing vagus trunks
# Rebuilds the File/Origin table: removes bones.exf, muscles.exf and the
# individual bone rows (scapula/clavicle/vertebral column), replacing the
# bone rows with a single "manubrium.exf" row, inserts "veins.exf" and
# "brainstem.exf" rows, and turns the nervesWithVagus.exf origin cell into
# a hyperlink.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out all the old data rows (keep header row 1 intact).
$ws.Range("A2:B15").ClearContents()

# Re-populate rows whose strings already existed in the workbook, in the
# same relative order as before, skipping the rows that are being removed
# (bones.exf, muscles.exf, Left/Right scapula.exf, left/right_clavicle.exf,
# vertebral column.exf).
$ws.Cells.Item(2,1).Value = "Geometry_Fitter_Wholebody.exf"
$ws.Cells.Item(2,2).Value = "https://sparc.science/datasets/file/307/6?path=files/primary/mapclient_workflow/Organs/Geometry_Fitter_Wholebody.exf"

$ws.Cells.Item(3,1).Value = "diaphragm.exf"
$ws.Cells.Item(3,2).Value = "https://sparc.science/datasets/file/307/6?path=files/primary/mapclient_workflow/Organs/diaphragm.exf"

$ws.Cells.Item(4,1).Value = "gastrointestinal_tract.exf"
$ws.Cells.Item(4,2).Value = "https://sparc.science/datasets/file/307/6?path=files/primary/mapclient_workflow/Organs/gastrointestinal_tract.exf"

$ws.Cells.Item(5,1).Value = "lung.exf"
$ws.Cells.Item(5,2).Value = "https://sparc.science/datasets/file/307/6?path=files/primary/mapclient_workflow/Organs/lung.exf"

$ws.Cells.Item(6,1).Value = "nervesWithVagus.exf"

$ws.Cells.Item(7,1).Value = "heart.exf"
$ws.Cells.Item(7,2).Value = "https://sparc.science/datasets/file/307/6?path=files/primary/mapclient_workflow/Organs/heart.exf"

$ws.Cells.Item(9,1).Value = "arteries.exf"
$ws.Cells.Item(9,2).Value = "https://sparc.science/datasets/file/307/6?path=files/primary/mapclient_workflow/Organs/arteries.exf"

# nervesWithVagus.exf's origin becomes a hyperlink to a (slightly altered,
# trailing-space-padded) URL.
$nervesUrl = "https://sparc.science/datasets/file/307/6?path=files/primary/mapclient_workflow/Organs/nervesWithVagus.exf  "
$ws.Cells.Item(6,2).Value = $nervesUrl

# New veins.exf row.
$ws.Cells.Item(8,1).Value = "veins.exf"
$ws.Cells.Item(8,2).Value = "https://sparc.science/datasets/file/307/6?path=files/primary/mapclient_workflow/Organs/veins.exf"

# New manubrium.exf row (entered before the brainstem.exf row below it).
$ws.Cells.Item(11,1).Value = "manubrium.exf"
$ws.Cells.Item(11,2).Value = "converted to exf from obj extracted from https://lifesciencedb.jp/bp3d/"

# New brainstem.exf row.
$ws.Cells.Item(10,1).Value = "brainstem.exf"
$ws.Cells.Item(10,2).Value = "https://sparc.science/datasets/file/307/6?path=files/primary/mapclient_workflow/Organs/brainstem.exf"

# Turn the nervesWithVagus.exf origin cell into an actual hyperlink
# (this also applies the built-in "Hyperlink" cell style).
$ws.Hyperlinks.Add($ws.Range("B6"), $nervesUrl)
